$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string value used for the "Type of Specialty" column on row 17
$ws.Range("G17").Value = "nmatx"

# Row 11: add Limit value in column H (between existing G11 "nforx" and I11)
$ws.Range("H11").Value = 6

# Row 17: update Markdown and Has Specialty, add new specialty columns (buy n get m at x%)
$ws.Range("E17").Value = 0.49
$ws.Range("F17").Value = $true
$ws.Range("H17").Value = 6
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 0.5

# Update the window view position/selection to reflect where the author was working
$ws.Range("H16").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
